$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 28 first (SC 92), then row 26 (RM 232), so row indices don't shift
# before the second deletion affects the first.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# Fill in the previously-missing F value for "SC 5", which is now row 26
$ws.Range("F26").Value = 17.38
